# Update the worksheet: refresh the date and every division-problem
# answer cell in the single practice table.

$d = $word.ActiveDocument

# --- Header date -----------------------------------------------------
$d.Content.Find.Execute("2025-06-20 Friday", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "2025-06-21 Saturday", 2)

# --- Table answer cells ------------------------------------------------
# The phrase "87÷2=43, 1" appears twice in the original table (row 1,
# col 1 and row 1, col 5) and the two occurrences are replaced with
# *different* values, so those two cells are addressed directly rather
# than via a global Find/Replace.
$t = $d.Tables.Item(1)
$t.Cell(1, 1).Range.Text = "37÷8=4, 5"
$t.Cell(1, 5).Range.Text = "90÷7=12, 6"

# All remaining answers are unique strings in the document, so a plain
# Find/Replace (restricted to whole-text matches) is unambiguous.
$replacements = @(
    @("12÷3=4, 0", "43÷4=10, 3"),
    @("86÷2=43, 0", "88÷7=12, 4"),
    @("33÷5=6, 3", "73÷4=18, 1"),
    @("33÷3=11, 0", "88÷2=44, 0"),
    @("99÷2=49, 1", "29÷8=3, 5"),
    @("90÷6=15, 0", "62÷7=8, 6"),
    @("47÷4=11, 3", "49÷8=6, 1"),
    @("19÷7=2, 5", "67÷9=7, 4"),
    @("74÷2=37, 0", "51÷9=5, 6"),
    @("75÷7=10, 5", "76÷9=8, 4"),
    @("58÷6=9, 4", "44÷5=8, 4"),
    @("66÷6=11, 0", "43÷5=8, 3"),
    @("44÷7=6, 2", "18÷5=3, 3"),
    @("42÷9=4, 6", "68÷3=22, 2"),
    @("93÷2=46, 1", "37÷6=6, 1"),
    @("85÷4=21, 1", "81÷9=9, 0"),
    @("82÷7=11, 5", "89÷6=14, 5"),
    @("88÷4=22, 0", "69÷6=11, 3"),
    @("67÷5=13, 2", "40÷5=8, 0"),
    @("76÷6=12, 4", "24÷8=3, 0"),
    @("66÷5=13, 1", "76÷7=10, 6"),
    @("30÷5=6, 0", "33÷2=16, 1"),
    @("39÷2=19, 1", "17÷4=4, 1")
)

foreach ($pair in $replacements) {
    $d.Content.Find.Execute($pair[0], $false, $false, $false, $false, $false, `
                             $true, 1, $false, $pair[1], 2)
}
